$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: relocate the trailing bookmark ("_GoBack") that currently sits
# between the "...Bait Store»" run and the ", именуемый в дальнейшем..." run
# of the introductory paragraph, and insert a brand-new paragraph
# ("Марка «Bait Store», ...") right after that introductory paragraph, so
# the bookmark ends up alone inside the (pre-existing) empty paragraph that
# used to immediately follow it.
# ---------------------------------------------------------------------------

# Locate the introductory paragraph by searching for its closing sentence,
# then expand the found range back out to the whole paragraph.
$introRange = $d.Content
$introRange.Find.ClearFormatting()
[void]$introRange.Find.Execute("способом.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
[void]$introRange.Expand(4)

# The paragraph right after it is the (pre-existing, empty) paragraph that
# currently hosts nothing but will end up hosting the relocated bookmark.
$introParaIndex = $d.Range(0, $introRange.Start).Paragraphs.Count + 1
$bookmarkParaIndex = $introParaIndex + 1

# Remove the old "_GoBack" bookmark from wherever it currently sits (between
# "Bait Store»" and the ", именуемый..." run, inside the introductory
# paragraph).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# Insert a brand-new paragraph *before* that pre-existing empty paragraph
# (rather than splitting the introductory paragraph itself) so the empty
# paragraph keeps being the very same, already-existing node -- inserting
# the bookmark into a paragraph freshly minted in this same script
# mis-places its end, so anchoring on the old node avoids that.
$anchor = $d.Paragraphs.Item($bookmarkParaIndex).Range
$anchor.Collapse(1)
$anchor.InsertParagraphBefore()

# The brand new paragraph now occupies the old empty paragraph's former
# index; fill it with the new sentence.
$newPara = $d.Paragraphs.Item($bookmarkParaIndex).Range
$newPara.InsertAfter("Марка «Bait Store», контактные данные и реквизиты интернет-магазина, условия работы с интернет-магазином, а также все товары и материалы, представленные на сайте, носят исключительно вымышленный характер и являются частью общего выдуманного контента (любые совпадения с реальными лицами, проектами и с прочими обстоятельствами случайны).")

# The original empty paragraph got pushed one slot further down; put the
# "_GoBack" bookmark back at its (now relocated) start.
$bookmarkHost = $d.Paragraphs.Item($bookmarkParaIndex + 1).Range
$bookmarkHost.Collapse(1)
$d.Bookmarks.Add("_GoBack", $bookmarkHost)

# ---------------------------------------------------------------------------
# Change 2: merge the two runs split by a <w:lastRenderedPageBreak/> into a
# single run/text ("... Производство спорадически транслирует анализ
# рыночных цен. Стимулирование ...").
# ---------------------------------------------------------------------------

$mergeRange = $d.Content
$mergeRange.Find.ClearFormatting()
[void]$mergeRange.Find.Execute( `
    ". Производство спорадически транслирует анализ рыночных цен. Стимулирование ", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    ". Производство спорадически транслирует анализ рыночных цен. Стимулирование ", 2)
